$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rail gun (row 22): improve images / stats tweak ---
$ws.Range("G22").Value = 80
$ws.Range("N22").Value = 55
$ws.Range("O22").Value = 72

# --- Balance plasma weapons ---
# Plasma Pistol (row 19)
$ws.Range("L19").Value = 12
$ws.Range("M19").Value = 16
$ws.Range("N19").Value = 25
$ws.Range("O19").Value = 35

# Plasma Rifle (row 20)
$ws.Range("L20").Value = 16
$ws.Range("M20").Value = 26
$ws.Range("N20").Value = 40
$ws.Range("O20").Value = 50

# Heavy Plasma (row 21)
$ws.Range("L21").Value = 18
$ws.Range("M21").Value = 28
$ws.Range("N21").Value = 40
$ws.Range("O21").Value = 50

# --- Balance sniper rifle (row 4) ---
$ws.Range("L4").Value = 9

# --- Balance laser weapons ---
# Laser pistol (row 5)
$ws.Range("F5").Value = 51
# Laser rifle (row 7)
$ws.Range("G7").Value = 99

# --- View / selection updates ---
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("M20").Select()
